# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Sheet "展览": F3 209->213, F4 827->832, F6 24->26
# Sheet "全部类型": F4 209->213, F5 827->832, F7 24->26

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 213
$wsExhibit.Range("F4").Value = 832
$wsExhibit.Range("F6").Value = 26

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 213
$wsAll.Range("F5").Value = 832
$wsAll.Range("F7").Value = 26
